# Appends new ticker rows to the end of the data range, as captured by the
# diff: rows 484-488 with new ticker symbols, updating the used dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTickers = @(
    "MNT-USD",
    "IMX-USD",
    "TAO-USD",
    "GRT-USD",
    "PEPE-USD"
)

$startRow = 484
for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
